$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.446.90'
$ws.Range('E2').Value = '  +1.32%  '

$ws.Range('D3').Value = '1.637.61'
$ws.Range('E3').Value = '  +2.38%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').Value = '''1.002'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.01%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '''306.47'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.19%  '

$ws.Range('E7').Value = '  -0.38%  '

$ws.Range('E8').Value = '  +0.54%  '

$ws.Range('D9').Value = '''0.3637'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.83%  '

$ws.Range('D10').Value = '''1.261'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.01%  '

$ws.Range('D11').Value = '''0.08147'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.48%  '

$ws.Range('D12').Value = '''1.002'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').Value = '''22.88'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.37%  '

$ws.Range('D14').Value = '''6.616'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.89%  '

$ws.Range('D15').Value = '''0.00001272'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.30%  '

$ws.Range('D16').Value = '''7.353'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.26%  '

$ws.Range('D17').Value = '1.640.07'
$ws.Range('E17').Value = '  +2.46%  '

$ws.Range('D18').Value = '''94.59'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.84%  '

$ws.Range('D19').Value = '''0.06947'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.01%  '

$ws.Range('D20').Value = '''18.12'
$ws.Range('D20').Style = "Normal"

$ws.Range('D21').Value = '''6.536'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.19%  '

$ws.Range('D22').Value = '''1.002'
$ws.Range('D22').Style = "Normal"

$ws.Range('D23').Value = '23.466.82'
$ws.Range('E23').Value = '  +1.37%  '

$ws.Range('D24').Value = '''12.75'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.28%  '

$ws.Range('D25').Value = '''3.100'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.14%  '

$ws.Range('E26').Value = '  +1.78%  '

$ws.Range('D27').Value = '''21.22'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.17%  '

$ws.Range('D28').Value = '''150.74'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.16%  '

$ws.Range('D29').Value = '''5.348'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.81%  '

$ws.Range('D30').Value = '''135.09'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.31%  '

$ws.Range('D31').Value = '''2.295'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.84%  '

$ws.Range('D32').Value = '1.823.49'
$ws.Range('E32').Value = '  +2.48%  '

$ws.Range('D33').Value = '''6.761'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.11%  '

$ws.Range('D34').Value = '''0.9626'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.04%  '

$ws.Range('D35').Value = '''0.02813'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.13%  '

$ws.Range('D36').Value = '''10.32'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.45%  '

$ws.Range('D37').Value = '''0.07312'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.27%  '

$ws.Range('D38').Value = '''0.2525'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.83%  '

$ws.Range('D39').Value = '''0.08835'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.43%  '

$ws.Range('E40').Value = '  +1.12%  '

$ws.Range('D41').Value = '''1.377'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.25%  '

$ws.Range('D42').Value = '''0.7089'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.28%  '

$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''12.49'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.34%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''16.13'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.15%  '

$ws.Range('D45').Value = '''0.6524'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.22%  '

$ws.Range('D46').Value = '''2.331'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.20%  '

$ws.Range('E47').Value = '  +0.05%  '

$ws.Range('D48').Value = '''4.024'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.33%  '

$ws.Range('D49').Value = '''0.07964'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.16%  '

$ws.Range('D50').Value = '''128.89'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.32%  '

$ws.Range('D51').Value = '''1.204'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.28%  '
